$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-09 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-10 Sunday", 2)
$d.Content.Find.Execute("241×5=1205", $true, $false, $false, $false, $false, $true, 1, $false, "264×2=528", 2)
$d.Content.Find.Execute("533×8=4264", $true, $false, $false, $false, $false, $true, 1, $false, "978×8=7824", 2)
$d.Content.Find.Execute("827×3=2481", $true, $false, $false, $false, $false, $true, 1, $false, "398×8=3184", 2)
$d.Content.Find.Execute("723×4=2892", $true, $false, $false, $false, $false, $true, 1, $false, "756×6=4536", 2)
$d.Content.Find.Execute("547×3=1641", $true, $false, $false, $false, $false, $true, 1, $false, "385×8=3080", 2)
$d.Content.Find.Execute("523×7=3661", $true, $false, $false, $false, $false, $true, 1, $false, "897×6=5382", 2)
$d.Content.Find.Execute("498×7=3486", $true, $false, $false, $false, $false, $true, 1, $false, "968×4=3872", 2)
$d.Content.Find.Execute("157×5=785", $true, $false, $false, $false, $false, $true, 1, $false, "198×9=1782", 2)
$d.Content.Find.Execute("660×5=3300", $true, $false, $false, $false, $false, $true, 1, $false, "679×4=2716", 2)
$d.Content.Find.Execute("192×4=768", $true, $false, $false, $false, $false, $true, 1, $false, "566×2=1132", 2)
$d.Content.Find.Execute("178×8=1424", $true, $false, $false, $false, $false, $true, 1, $false, "529×9=4761", 2)
$d.Content.Find.Execute("479×9=4311", $true, $false, $false, $false, $false, $true, 1, $false, "434×8=3472", 2)
$d.Content.Find.Execute("446×8=3568", $true, $false, $false, $false, $false, $true, 1, $false, "423×4=1692", 2)
$d.Content.Find.Execute("219×4=876", $true, $false, $false, $false, $false, $true, 1, $false, "183×9=1647", 2)
$d.Content.Find.Execute("743×4=2972", $true, $false, $false, $false, $false, $true, 1, $false, "448×5=2240", 2)
$d.Content.Find.Execute("759×6=4554", $true, $false, $false, $false, $false, $true, 1, $false, "102×9=918", 2)
$d.Content.Find.Execute("472×2=944", $true, $false, $false, $false, $false, $true, 1, $false, "143×6=858", 2)
$d.Content.Find.Execute("588×8=4704", $true, $false, $false, $false, $false, $true, 1, $false, "272×7=1904", 2)
$d.Content.Find.Execute("517×7=3619", $true, $false, $false, $false, $false, $true, 1, $false, "664×5=3320", 2)
$d.Content.Find.Execute("624×9=5616", $true, $false, $false, $false, $false, $true, 1, $false, "299×5=1495", 2)
$d.Content.Find.Execute("712×9=6408", $true, $false, $false, $false, $false, $true, 1, $false, "926×6=5556", 2)
$d.Content.Find.Execute("746×8=5968", $true, $false, $false, $false, $false, $true, 1, $false, "138×4=552", 2)
$d.Content.Find.Execute("270×9=2430", $true, $false, $false, $false, $false, $true, 1, $false, "163×3=489", 2)
$d.Content.Find.Execute("612×9=5508", $true, $false, $false, $false, $false, $true, 1, $false, "423×7=2961", 2)
$d.Content.Find.Execute("511×8=4088", $true, $false, $false, $false, $false, $true, 1, $false, "813×7=5691", 2)
